$d = $word.ActiveDocument

# 1. Update the GDK version text in the compatibility paragraph.
#    "June 2020" -> "October 2021" (this also splits the run in the
#    real document because the editor replaced only the selected text).
$d.Content.Find.Execute("Microsoft Game Development Kit (June 2020)", $true, $false, $false, $false, $false, $true, 1, $false, "Microsoft Game Development Kit (October 2021)", 2)

# 2. Update the copyright year shown in both footers (cached field result
#    of the `DATE \@ "yyyy"` field) from 2021 to 2022.
foreach ($sec in $d.Sections) {
    foreach ($idx in 1, 2, 3) {
        $ftr = $sec.Footers.Item($idx)
        if ($ftr -and $ftr.Exists) {
            $ftr.Range.Find.Execute("2021", $true, $false, $false, $false, $false, $true, 1, $false, "2022", 2)
        }
    }
}
